# Generate Report for Handback
# Updates the localization-status workbook after a handback transform
# failure: the "Ready for handoff" status for the 2652da41... file becomes
# "Handback transform failed", and the per-locale "Error Detail" column is
# populated with the mismatch explanation (and widened to fit).

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$zhError = "Handback file name: censubje.k5v is different with handoff file name: 2652da41-441c-4de1-a54b-3b9db6c9be5d.50d29f61f113ef8ac4dd6a6001406ac37f30765a.zh-cn."
$deError  = "Handback file name: censubje.k5v is different with handoff file name: 2652da41-441c-4de1-a54b-3b9db6c9be5d.50d29f61f113ef8ac4dd6a6001406ac37f30765a.de-de."

# --- Overview sheet: update the zh-cn / de-de status for the 2652da41 row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 39.17
